$d = $word.ActiveDocument

$replacements = @(
    @{old = "23×46=1058"; new = "67×62=4154"},
    @{old = "44×85=3740"; new = "78×63=4914"},
    @{old = "49×15=735";  new = "38×66=2508"},
    @{old = "48×48=2304"; new = "63×82=5166"},
    @{old = "60×62=3720"; new = "90×76=6840"},
    @{old = "41×24=984";  new = "77×38=2926"},
    @{old = "34×36=1224"; new = "75×67=5025"},
    @{old = "43×82=3526"; new = "89×13=1157"},
    @{old = "78×91=7098"; new = "56×99=5544"},
    @{old = "23×52=1196"; new = "64×75=4800"},
    @{old = "47×61=2867"; new = "14×55=770"},
    @{old = "19×83=1577"; new = "73×17=1241"},
    @{old = "39×87=3393"; new = "66×52=3432"},
    @{old = "58×75=4350"; new = "12×64=768"},
    @{old = "36×13=468";  new = "73×89=6497"},
    @{old = "28×85=2380"; new = "96×73=7008"},
    @{old = "19×42=798";  new = "65×98=6370"},
    @{old = "37×66=2442"; new = "19×50=950"},
    @{old = "11×87=957";  new = "51×92=4692"},
    @{old = "80×86=6880"; new = "24×87=2088"},
    @{old = "48×13=624";  new = "79×37=2923"},
    @{old = "68×14=952";  new = "78×48=3744"},
    @{old = "56×20=1120"; new = "75×35=2625"},
    @{old = "68×60=4080"; new = "29×61=1769"},
    @{old = "47×35=1645"; new = "82×79=6478"}
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $r.new, 2)
}
